# Rename speaker transcripts in column D (Speaker) on the active sheet.
# Full names are replaced by their short codes:
#   ANTOINETTE VILLARIN -> T
#   PATTY FERRANT        -> T2
#   STUDENT              -> S
#   STUDENTS             -> SS

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "ANTOINETTE VILLARIN" = "T"
    "PATTY FERRANT" = "T2"
    "STUDENT" = "S"
    "STUDENTS" = "SS"
}

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $cur = $cell.Value2
    if ($cur -ne $null -and $map.ContainsKey($cur)) {
        $cell.Value() = $map[$cur]
    }
}
